# account_bank_statement_import_adyen test file migration (12.0 MIG)
#  - Gross/Net currency sample data changed from EUR to USD
#  - Gross Debit (GC) sample value updated on row 10
#  - Creation Date format normalised to lower-case date/time tokens
#  - Active cell selection moved
#  - Keep gridlines visible (matches the canonical sheet view)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Currency sample data: EUR -> USD -------------------------------------
# Column K = "Gross Currency", Column O = "Net Currency". All cells that
# currently read "EUR" (K5:K17, K19:K23, O5:O25 - K18 is "GBP" and stays as is)
# become "USD".
$ws.Range("K5:K17").Value = "USD"
$ws.Range("K19:K23").Value = "USD"
$ws.Range("O5:O25").Value = "USD"

# --- Gross Debit (GC) value on row 10: 666 -> 1598 ------------------------
$ws.Range("M10").Value = 1598

# --- Creation Date number format: upper-case tokens -> lower-case --------
$ws.Range("G5:G25").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# --- View state: keep gridlines on, move the active selection ------------
$excel.ActiveWindow.DisplayGridlines = $True
[void]$ws.Range("L9").Select()
